$wb = $excel.ActiveWorkbook

# 1. "Files" sheet: a new row is added listing the bag-info.txt file that is now
#    included in the generated bag (dimension grows from A1:D1 to A1:D2).
$wsFiles = $wb.Worksheets.Item("Files")
$wsFiles.Range("A2").Value = "bag-info.txt"

# 2. "Licenses" sheet: the CC BY-NC-SA 3.0 AU description text gets extra blank
#    (CR) lines inserted before each paragraph, and a new "TYPE:" value
#    ("ScholarlyWork") is recorded in D2.
$wsLicenses = $wb.Worksheets.Item("Licenses")
$licenseText = @"
Attribution-NonCommercial-ShareAlike 3.0 Australia (CC BY-NC-SA 3.0 AU)_x000d_
_x000d_
_x000d_
_x000d_
_x000d_
_x000d_
_x000d_
_x000d_
_x000d_
_x000d_
_x000d_
_x000d_
_x000d_
_x000d_
_x000d_
_x000d_
This is a human-readable summary of (and not a substitute for) the license. Disclaimer._x000d_
_x000d_
_x000d_
_x000d_
_x000d_
_x000d_
_x000d_
_x000d_
_x000d_
_x000d_
_x000d_
_x000d_
_x000d_
_x000d_
_x000d_
_x000d_
_x000d_
_x000d_
_x000d_
_x000d_
_x000d_
_x000d_
_x000d_
_x000d_
_x000d_
_x000d_
_x000d_
_x000d_
_x000d_
_x000d_
_x000d_
_x000d_
You are free to:_x000d_
_x000d_
_x000d_
_x000d_
_x000d_
_x000d_
_x000d_
_x000d_
_x000d_
_x000d_
_x000d_
_x000d_
_x000d_
_x000d_
_x000d_
_x000d_
_x000d_
_x000d_
_x000d_
_x000d_
_x000d_
_x000d_
_x000d_
_x000d_
_x000d_
_x000d_
_x000d_
_x000d_
_x000d_
_x000d_
_x000d_
_x000d_
Share — copy and redistribute the material in any medium or format_x000d_
_x000d_
_x000d_
_x000d_
_x000d_
_x000d_
_x000d_
_x000d_
_x000d_
_x000d_
_x000d_
_x000d_
_x000d_
_x000d_
_x000d_
_x000d_
Adapt — remix, transform, and build upon the material_x000d_
_x000d_
_x000d_
_x000d_
_x000d_
_x000d_
_x000d_
_x000d_
_x000d_
_x000d_
_x000d_
_x000d_
_x000d_
_x000d_
_x000d_
_x000d_
The licensor cannot revoke these freedoms as long as you follow the license terms._x000d_
_x000d_
_x000d_
_x000d_
_x000d_
_x000d_
_x000d_
_x000d_
_x000d_
_x000d_
_x000d_
_x000d_
_x000d_
_x000d_
_x000d_
_x000d_
Under the following terms:_x000d_
_x000d_
_x000d_
_x000d_
_x000d_
_x000d_
_x000d_
_x000d_
_x000d_
_x000d_
_x000d_
_x000d_
_x000d_
_x000d_
_x000d_
_x000d_
_x000d_
_x000d_
_x000d_
_x000d_
_x000d_
_x000d_
_x000d_
_x000d_
_x000d_
_x000d_
_x000d_
_x000d_
_x000d_
_x000d_
_x000d_
_x000d_
Attribution — You must give appropriate credit, provide a link to the license, and indicate if changes were made. You may do so in any reasonable manner, but not in any way that suggests the licensor endorses you or your use._x000d_
_x000d_
_x000d_
_x000d_
_x000d_
_x000d_
_x000d_
_x000d_
_x000d_
_x000d_
_x000d_
_x000d_
_x000d_
_x000d_
_x000d_
_x000d_
Non-Commercial — You may not use the material for commercial purposes._x000d_
_x000d_
_x000d_
_x000d_
_x000d_
_x000d_
_x000d_
_x000d_
_x000d_
_x000d_
_x000d_
_x000d_
_x000d_
_x000d_
_x000d_
_x000d_
ShareAlike — If you remix, transform, or build upon the material, you must distribute your contributions under the same license as the original._x000d_
_x000d_
_x000d_
_x000d_
_x000d_
_x000d_
_x000d_
_x000d_
_x000d_
_x000d_
_x000d_
_x000d_
_x000d_
_x000d_
_x000d_
_x000d_
No additional restrictions — You may not apply legal terms or technological measures that legally restrict others from doing anything the license permits.
"@
$wsLicenses.Range("C2").Value = $licenseText
$wsLicenses.Range("D2").Value = "ScholarlyWork"

# 3. "Projects" sheet: trim the trailing space from E3 ("Project " -> "Project")
$wsProjects = $wb.Worksheets.Item("Projects")
$wsProjects.Range("E3").Value = "Project"
